$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("90÷9=10, 0", "42÷5=8, 2", "34÷8=4, 2", "33÷8=4, 1", "96÷7=13, 5")
    5  = @("36÷3=12, 0", "34÷4=8, 2", "94÷5=18, 4", "35÷6=5, 5", "34÷4=8, 2")
    9  = @("37÷4=9, 1", "53÷5=10, 3", "95÷4=23, 3", "12÷2=6, 0", "77÷4=19, 1")
    13 = @("13÷6=2, 1", "39÷3=13, 0", "57÷6=9, 3", "60÷4=15, 0", "27÷7=3, 6")
    17 = @("86÷4=21, 2", "34÷9=3, 7", "69÷2=34, 1", "88÷5=17, 3", "73÷3=24, 1")
}

foreach ($rowIndex in $newValues.Keys) {
    $values = $newValues[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
